# FAST_holdings.xlsx update
# - Refresh the "as of" date in the confidential disclosure note (A13)
# - Refresh Weight (col D) and Percent Change (col E) figures for rows 2-10
#
# The sheet ships with sheetProtection, so locked cells (the default for
# every cell here) cannot be written to until the sheet is unprotected.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Confidential notice: bump the model-holdings-as-of date one day forward.
$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-19 for illustrative purposes only and are subject to change."

# Row 2 (ARKG)
$ws.Range("D2").Value = 0.09202996156374289
$ws.Range("E2").Value = -0.0183568677792042

# Row 3 (MTUM)
$ws.Range("D3").Value = 0.1051656416260982
$ws.Range("E3").Value = 0.0002468221646303448

# Row 4 (QUAL)
$ws.Range("D4").Value = 0.1198349602567391
$ws.Range("E4").Value = -0.001424614166996352

# Row 5 (SIZE)
$ws.Range("D5").Value = 0.1415014112594686
$ws.Range("E5").Value = -0.00543739005277466

# Row 6 (USMV)
$ws.Range("D6").Value = 0.1378302961751589
$ws.Range("E6").Value = -0.004164931278633932

# Row 7 (VLUE)
$ws.Range("D7").Value = 0.1470795494132282
$ws.Range("E7").Value = -0.002850898032880189

# Row 8 (SNSR)
$ws.Range("D8").Value = 0.1267893297871168
$ws.Range("E8").Value = -0.004829459704195682

# Row 9 (FITE)
$ws.Range("D9").Value = 0.1297688499184475
$ws.Range("E9").Value = -0.004846065553383183

# Row 10 (Total) - only Percent Change changes, Weight stays 1
$ws.Range("E10").Value = -0.004838096395362035

# Restore sheet protection.
$ws.Protect()
